$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos price/volume data (and a few reordered coin rows) per latest scrape
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.405.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.696.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.684.02"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.74%  "
$ws.Range("E8").Value = "  +3.91%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.612"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000286"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.292.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "680.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.701.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.523.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.68%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.941"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.06%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "575.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.726.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.146"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.65%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0769"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0469"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.56%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.24%  "
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.50%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
